$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly log rolled forward by one week: a new week (44664) was appended
# and the oldest week (44636) pushed out, shifting each row's data up by one
# record while keeping row 2 at the top of the table.
#
# New row data (row number -> Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad de comercializacion,
# Precio $/Kg, Kg/unidad)
$rows = @(
    @{ Row = 2; D = 44664; L = "Segunda"; M = 150; N = 29000; O = 30000; P = 29500; Q = "`$/caja 18 kilos"; S = 1639; T = 18 },
    @{ Row = 3; D = 44679; L = "Segunda"; M = 200; N = 29000; O = 30000; P = 29500; Q = "`$/caja 20 kilos"; S = 1475; T = 20 },
    @{ Row = 4; D = 44679; L = "Tercera"; M = 200; N = 24000; O = 25000; P = 24500; Q = "`$/caja 20 kilos"; S = 1225; T = 20 },
    @{ Row = 5; D = 44650; L = "Primera"; M = 160; N = 31000; O = 32000; P = 31500; Q = "`$/caja 20 kilos"; S = 1575; T = 20 },
    @{ Row = 6; D = 44650; L = "Segunda"; M = 250; N = 29000; O = 30000; P = 29500; Q = "`$/caja 20 kilos"; S = 1475; T = 20 },
    @{ Row = 7; D = 44643; L = "Primera"; M = 160; N = 28000; O = 30000; P = 29000; Q = "`$/caja 20 kilos"; S = 1450; T = 20 },
    @{ Row = 8; D = 44671; L = "Segunda"; M = 200; N = 29000; O = 30000; P = 29500; Q = "`$/caja 20 kilos"; S = 1475; T = 20 },
    @{ Row = 9; D = 44636; L = "Primera"; M = 200; N = 29000; O = 30000; P = 29500; Q = "`$/caja 20 kilos"; S = 1475; T = 20 }
)

foreach ($r in $rows) {
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("L" + $r.Row).Value = $r.L
    $ws.Range("M" + $r.Row).Value = $r.M
    $ws.Range("N" + $r.Row).Value = $r.N
    $ws.Range("O" + $r.Row).Value = $r.O
    $ws.Range("P" + $r.Row).Value = $r.P
    $ws.Range("Q" + $r.Row).Value = $r.Q
    $ws.Range("S" + $r.Row).Value = $r.S
    $ws.Range("T" + $r.Row).Value = $r.T
}
